$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.162972569465637
$ws.Range("B1").Value = 2.389852523803711
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.783748269081116
$ws.Range("E1").Value = 1.195452213287354
